# Regression test fix: swap the bus-name (column A) and angle (column E)
# values between paired rows on the "Bus" sheet so that the phase-A and
# phase-C entries of each bus are correctly ordered.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bus")

# Row pairs (1-based worksheet rows) whose A/E values must be swapped.
$pairs = @(
    @(3, 5),
    @(6, 8),
    @(9, 11),
    @(12, 13),
    @(14, 15),
    @(16, 18),
    @(20, 22),
    @(23, 25),
    @(26, 28),
    @(29, 31),
    @(32, 33),
    @(34, 36),
    @(37, 39),
    @(42, 44),
    @(45, 47),
    @(48, 50),
    @(52, 54)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $aCell1 = $ws.Range("A$r1")
    $eCell1 = $ws.Range("E$r1")
    $aCell2 = $ws.Range("A$r2")
    $eCell2 = $ws.Range("E$r2")

    $a1 = $aCell1.Value()
    $e1 = $eCell1.Value()
    $a2 = $aCell2.Value()
    $e2 = $eCell2.Value()

    $aCell1.Value = $a2
    $eCell1.Value = $e2
    $aCell2.Value = $a1
    $eCell2.Value = $e1
}
